$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of 1-based row index -> new cell text (single column table).
$rowUpdates = [ordered]@{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "203"
    6  = "0.00009"
    7  = "0.00004"
    8  = "0.00000"
    9  = "0.00009"
    10 = "0.00009"
    11 = "0.00009"
    12 = "0.00693"
    44 = "100"
    45 = "0.01"
    46 = "371"
}

foreach ($rowIndex in $rowUpdates.Keys) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $rowUpdates[$rowIndex]
}
